$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated price (D) and 1h-volume-change (E) figures scraped for this run.
# D-column values are forced to Text format before assignment (and reset back to General/Normal
# afterwards) so that numeric-looking strings (e.g. "237.62") are preserved verbatim as text
# instead of being auto-converted to Excel numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '91.099.60'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.23%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.142.86'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.81%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.62'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +8.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '636.90'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.27%  '
$ws.Range("E7").Value = '  +6.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.368'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.74%  '
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.143.92'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.95%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.726'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.10%  '
$ws.Range("E12").Value = '  +2.17%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '36.53'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.41%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000250'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.03%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.58'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.69%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '90.998.08'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.726.84'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.99%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.128.80'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.88%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.75'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.33%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000215'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.39'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.96%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '445.95'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.01%  '
$ws.Range("E23").Value = '  +8.55%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.02'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.43%  '
$ws.Range("E25").Value = '  -3.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '90.59'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.48'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.26%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.310.41'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.89%  '
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.73'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.160'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.998'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +11.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.201'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +29.37%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.00'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +12.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.85'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '515.85'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.77%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.150'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.18'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.25%  '
$ws.Range("E39").Value = '  +4.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.31'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.41%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.420'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.34%  '
$ws.Range("E42").Value = '  -0.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0859'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.61%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.39'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +48.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.95'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.92%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.701'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +11.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '150.98'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.63%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '45.72'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.95%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.58'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.78%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.35'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.86%  '
